$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-19:
# all these cells currently hold serial date 45207 (2023-10-08) and
# should be bumped to 45208 (2023-10-09).
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
